$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 3972.90162440078
$ws.Range("C2").Value = 5569.6059699143
$ws.Range("F2").Value = 24.7318882713966
$ws.Range("B3").Value = 4035.78924951556
$ws.Range("C3").Value = 5174.80556047602
$ws.Range("F3").Value = 107.787898540019
$ws.Range("B4").Value = 862.43729410157
$ws.Range("C4").Value = 3093.59811189519
$ws.Range("F4").Value = 11.5033965747341
$ws.Range("B5").Value = 746.816712468479
$ws.Range("C5").Value = 2768.74755695588
$ws.Range("F5").Value = -2.9140998130249
$ws.Range("B6").Value = 4021.33103157071
$ws.Range("C6").Value = 5371.08183364754
$ws.Range("F6").Value = 115.662013378201
$ws.Range("B7").Value = 4015.12679520747
$ws.Range("C7").Value = 5070.80726511363
$ws.Range("F7").Value = 155.170402454423
$ws.Range("B8").Value = 4355.11679040039
$ws.Range("C8").Value = 5090.12048693026
$ws.Range("F8").Value = 141.808870230411
$ws.Range("B9").Value = 4355.11679040039
$ws.Range("C9").Value = 4900.99600552337
$ws.Range("F9").Value = 133.928683505124
$ws.Range("B10").Value = 5445.71283757669
$ws.Range("C10").Value = 4233.7273372665
$ws.Range("F10").Value = 142.288045695409
$ws.Range("B11").Value = 1503.61954698234
$ws.Range("C11").Value = 2167.19216528301
$ws.Range("F11").Value = 39.8173470125281
$ws.Range("B12").Value = 1400.00991904533
$ws.Range("C12").Value = 2209.22675481346
$ws.Range("F12").Value = 41.1664603236721
$ws.Range("B13").Value = 5841.30158853682
$ws.Range("C13").Value = 5086.84982201133
$ws.Range("F13").Value = 192.865878519771
$ws.Range("B14").Value = 5841.30158853682
$ws.Range("C14").Value = 5503.47865688923
$ws.Range("F14").Value = 210.22541330635
$ws.Range("B15").Value = 5841.30158853682
$ws.Range("C15").Value = 5681.48593630429
$ws.Range("F15").Value = 218.517799365311
